$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value in B1 from 3 to 2
$ws.Range("B1").Value = 2

# Move the active selection to A2
$ws.Range("A2").Select()
